$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column cells to keep their text representation
# (these values look numeric to Excel and would otherwise be
# auto-converted/normalized to a float, losing formatting such as
# trailing zeros or thousand-separator dots).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '95.096.53'
$ws.Range("E2").Value = '  +1.42%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.616.15'
$ws.Range("E3").Value = '  +3.89%  '

$ws.Range("E4").Value = '  -0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '236.66'
$ws.Range("E5").Value = '  +0.08%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '658.39'
$ws.Range("E6").Value = '  +5.20%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.46'
$ws.Range("E7").Value = '  +0.55%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.405'
$ws.Range("E8").Value = '  +2.24%  '

$ws.Range("E9").Value = '  -0.01%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.994'
$ws.Range("E10").Value = '  -2.55%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '3.613.91'
$ws.Range("E11").Value = '  +3.92%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '42.46'
$ws.Range("E12").Value = '  -3.98%  '

$ws.Range("E13").Value = '  +0.15%  '

$ws.Range("E14").Value = '  -0.24%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.305.64'
$ws.Range("E15").Value = '  +4.12%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '95.058.04'
$ws.Range("E16").Value = '  +1.51%  '

$ws.Range("E17").Value = '  +1.26%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.614.04'
$ws.Range("E18").Value = '  +3.82%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.92'
$ws.Range("E19").Value = '  -5.91%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.82'
$ws.Range("E20").Value = '  +8.37%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.99'
$ws.Range("E21").Value = '  -2.37%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.57'
$ws.Range("E22").Value = '  +4.77%  '

$ws.Range("E23").Value = '  -8.75%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '504.99'
$ws.Range("E24").Value = '  +0.36%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0000197'
$ws.Range("E25").Value = '  +5.61%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.61'
$ws.Range("E26").Value = '  -3.42%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '95.37'
$ws.Range("E27").Value = '  -0.30%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '3.810.48'
$ws.Range("E28").Value = '  +4.57%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '12.55'
$ws.Range("E29").Value = '  +2.68%  '

$ws.Range("E30").Value = '  +13.14%  '

$ws.Range("B31").Value = 'InternetComputer(DFINITY)'
$ws.Range("C31").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '11.27'
$ws.Range("E31").Value = '  -2.11%  '

$ws.Range("B32").Value = 'Dai'
$ws.Range("C32").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.00'
$ws.Range("E32").Value = '  +0.01%  '

$ws.Range("E33").Value = '  -3.50%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.996'
$ws.Range("E34").Value = '  +0.65%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '32.14'
$ws.Range("E35").Value = '  +8.16%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.177'
$ws.Range("E36").Value = '  -2.09%  '

$ws.Range("E37").Value = '  -0.64%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '570.26'
$ws.Range("E38").Value = '  -0.14%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '8.11'
$ws.Range("E39").Value = '  +6.26%  '

$ws.Range("E40").Value = '  +3.10%  '

$ws.Range("E41").Value = '  +0.00%  '

$ws.Range("B42").Value = 'Kaspa'
$ws.Range("C42").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.149'
$ws.Range("E42").Value = '  -0.55%  '

$ws.Range("B43").Value = 'ARBITRUM'
$ws.Range("C43").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.915'
$ws.Range("E43").Value = '  -0.92%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '35.27'
$ws.Range("E44").Value = '  +40.97%  '

$ws.Range("E45").Value = '  +0.21%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '23.69'
$ws.Range("E46").Value = '  -0.23%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '5.61'
$ws.Range("E47").Value = '  +0.57%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.23'
$ws.Range("E48").Value = '  +4.64%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0413'
$ws.Range("E49").Value = '  -3.27%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '3.55'
$ws.Range("E50").Value = '  -2.35%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '53.41'
$ws.Range("E51").Value = '  +0.27%  '

